$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (crearPartido): clear B6 first so shared string slot for
#     "No es necesario" becomes free and gets reused/renamed in place ---
$ws.Range("B6").Value = ""

# Apply the "black/white note" style (same as used in row 11) to the
# whole of rows 6, 9 and 26 by copying formats from the existing row 11.
$ws.Range("A11:D11").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Range("A9:D9").PasteSpecial(-4122)
$ws.Range("A26:D26").PasteSpecial(-4122)

# Re-apply the correct style to A2:D2 (must stay "Incorrecto") - PasteSpecial
# above did not touch row 2, so nothing to restore there.

# --- Fill in the new comments (column D) ---
$ws.Range("D6").Value = "Por ahora no se realizará, se hará directamente seleccionar equipos"
$ws.Range("D9").Value = "Por ahora no se realizará, se hará directamente seleccionar equipos"
$ws.Range("D26").Value = "Se llamará Create  Match en la parte front y estará fusionada con configuracion equipo"
$ws.Range("D2").Value = "Se llamara Create Team, solo corresponderá a la parte de seleccionar equipos"

# --- Column D width ---
$ws.Columns.Item(4).ColumnWidth = 77.1

# --- Selection ---
$ws.Range("E10").Select()
